# Adds two new trailing columns (AX: "Max Bootstrapped Demos",
# AY: "Number of Candidate Programs") to the evaluation log sheet, leaves
# them blank for the existing rows (2-11), and appends a new results row
# (12) for a "qwen2:7b-instruct-q5_K_M" / "llama3:70b" bootstrap run that
# populates all columns, including the two new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells, formatted like the rest of row 1 -------------------
$ws.Range("AW1").Copy() | Out-Null
$ws.Range("AX1:AY1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("AX1").Value = "Max Bootstrapped Demos"
$ws.Range("AY1").Value = "Number of Candidate Programs"

# --- Materialize (but leave blank) AX/AY for the pre-existing rows --------
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 50).Font.Bold = $false   # AX{r}
    $ws.Cells.Item($r, 51).Font.Bold = $false   # AY{r}
}

# --- New row 12: qwen2:7b-instruct-q5_K_M vs llama3:70b (bootstrap) -------
$ws.Range("A12").Value = "qwen2:7b-instruct-q5_K_M"
$ws.Range("B12").Value = "llama3:70b"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 200
$ws.Range("E12").Value = 2138.53
$ws.Range("F12").Value = 50.6
$ws.Range("G12").Value = 42.5
$ws.Range("H12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_match.txt'
$ws.Range("I12").Value = 50.6
$ws.Range("J12").Value = 90
$ws.Range("K12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_correct.txt'
$ws.Range("L12").Value = 31.33333333333333
$ws.Range("M12").Value = 94.84999999999999
$ws.Range("N12").Value = 51.25
$ws.Range("O12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_match.txt'
$ws.Range("P12").Value = 94.84999999999999
$ws.Range("Q12").Value = 91.25
$ws.Range("R12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_correct.txt'
$ws.Range("S12").Value = 42.33333333333334
$ws.Range("T12").Value = 0
$ws.Range("U12").Value = 92.77
$ws.Range("V12").Value = 57.5
$ws.Range("W12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_fewshot_match.txt'
$ws.Range("X12").Value = 92.77
$ws.Range("Y12").Value = 90
$ws.Range("Z12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_fewshot_correct.txt'
$ws.Range("AA12").Value = 40.66666666666666
$ws.Range("AB12").Value = 179.56
$ws.Range("AC12").Value = 43.75
$ws.Range("AD12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_fewshot_match.txt'
$ws.Range("AE12").Value = 179.56
$ws.Range("AF12").Value = 88.75
$ws.Range("AG12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_fewshot_correct.txt'
$ws.Range("AH12").Value = 31.33333333333333
$ws.Range("AI12").Value = 1138.61
$ws.Range("AJ12").Value = 199.39
$ws.Range("AK12").Value = 47.5
$ws.Range("AL12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_bootstrap_match.txt'
$ws.Range("AM12").Value = 199.39
$ws.Range("AN12").Value = 82.5
$ws.Range("AO12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_bootstrap_correct.txt'
$ws.Range("AP12").Value = 31.33333333333333
$ws.Range("AQ12").Value = 382.76
$ws.Range("AR12").Value = 46.25
$ws.Range("AS12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_bootstrap_match.txt'
$ws.Range("AT12").Value = 382.76
$ws.Range("AU12").Value = 85
$ws.Range("AV12").Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_bootstrap_correct.txt'
$ws.Range("AW12").Value = 31
$ws.Range("AX12").Value = 2
$ws.Range("AY12").Value = 2

Write-Output "edit complete"
